# Update the "as_of_utc" timestamp column (AA) on the data sheets
# from "2025-12-12 03:03:16" to "2025-12-12 07:02:38".

$wb = $excel.ActiveWorkbook

$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = "2025-12-12 07:02:38"
    }
}
